$d = $word.ActiveDocument

# --- Problem 1 section: collapse runs that were only split apart by now-removed
# proofing (spell-check) marks around "c&p" / "p&s" style tokens. The wording
# itself is unchanged; re-asserting the same text merges the runs and drops the
# spell-check annotations, same as Word does when it re-flows the paragraph. ---

$d.Content.Find.Execute("In order for the solution to work c&p or p&s can", $true, $false, $false, $false, $false, $true, 1, $false, "In order for the solution to work c&p or p&s can", 2)

$d.Content.Find.Execute("Man takes c leaving p&s together = Unsuccessful", $true, $false, $false, $false, $false, $true, 1, $false, "Man takes c leaving p&s together = Unsuccessful", 2)

$d.Content.Find.Execute("Man takes p leaving c&s together = Successful", $true, $false, $false, $false, $false, $true, 1, $false, "Man takes p leaving c&s together = Successful", 2)

$d.Content.Find.Execute("Man leaves c&p together and going to s = Unsuccessful", $true, $false, $false, $false, $false, $true, 1, $false, "Man leaves c&p together and going to s = Unsuccessful", 2)

$d.Content.Find.Execute("Man take s leaving c&p together = Unsuccessful", $true, $false, $false, $false, $false, $true, 1, $false, "Man take s leaving c&p together = Unsuccessful", 2)

# --- Problem 2 section: rewrite "...but to resolve is to break down when it
# will occur." into "...but to resolve this problem is to break down *when*
# it will occur." (note the single trailing space and the italic "when"). ---

$d.Content.Find.Execute("but to resolve is to break down", $true, $false, $false, $false, $false, $true, 1, $false, "but to resolve this problem is to break down", 2)

$d.Content.Find.Execute("when it will occur.  ", $true, $false, $false, $false, $false, $true, 1, $false, "when it will occur. ", 2)

$rng = $d.Content
$rng.Find.Execute("when it will occur", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$whenRng = $d.Range($rng.Start, $rng.Start + 4)
$whenRng.Italic = 1
